$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62: fill in OHLCV + change values (C62:H62) ---
$ws.Range("C62").Value = 591.55999999999995
$ws.Range("D62").Value = 592.77
$ws.Range("E62").Value = 586.99
$ws.Range("F62").Value = 587.73
$ws.Range("G62").Value = 66298636
$ws.Range("H62").Value = 0.15988055714988222

# --- Row 63: populate full data row (A63:II63) ---
$ws.Range("A63").Value = 45805
$ws.Range("B63").Value = 45806
$ws.Range("I63").Value = 0.19309999999999999
$ws.Range("J63").Value = 4.7
$ws.Range("K63").Value = 600
$ws.Range("L63").Value = 325951800
$ws.Range("M63").Value = 16594
$ws.Range("N63").Value = 331
$ws.Range("O63").Value = 16925
$ws.Range("P63").Value = 0.10397850693860272
$ws.Range("Q63").Value = 0.053536478369440504
$ws.Range("R63").Value = 45807
$ws.Range("S63").Value = 0.14306539785381656
$ws.Range("T63").Value = 45814
$ws.Range("U63").Value = 0.059495849362219071
$ws.Range("V63").Value = 45828
$ws.Range("W63").Value = 0.36221232368225686
$ws.Range("X63").Value = 11.333333333333334
$ws.Range("Y63").Value = 595
$ws.Range("Z63").Value = 324281545
$ws.Range("AA63").Value = 12427
$ws.Range("AB63").Value = 1102
$ws.Range("AC63").Value = 13529
$ws.Range("AD63").Value = 0.10344569619447817
$ws.Range("AE63").Value = 0.16051600338455865
$ws.Range("AF63").Value = 45806
$ws.Range("AG63").Value = 0.16051600338455865
$ws.Range("AH63").Value = 45807
$ws.Range("AI63").Value = 0.4472454766192861
$ws.Range("AJ63").Value = 45828
$ws.Range("AK63").Value = 0.18321050914918058
$ws.Range("AL63").Value = 8.6666666666666661
$ws.Range("AM63").Value = 592
$ws.Range("AN63").Value = 223910976
$ws.Range("AO63").Value = 11690
$ws.Range("AP63").Value = 633
$ws.Range("AQ63").Value = 12323
$ws.Range("AR63").Value = 0.071427520791863416
$ws.Range("AS63").Value = 0.62460398125883143
$ws.Range("AT63").Value = 45806
$ws.Range("AU63").Value = 0.62460398125883143
$ws.Range("AV63").Value = 45807
$ws.Range("AW63").Value = 0.19863407620416967
$ws.Range("AX63").Value = 45810
$ws.Range("AY63").Value = 0.05729938768933291
$ws.Range("AZ63").Value = 2.6666666666666665
$ws.Range("BA63").Value = 590
$ws.Range("BB63").Value = 192178340
$ws.Range("BC63").Value = 6119
$ws.Range("BD63").Value = -1631
$ws.Range("BE63").Value = 7750
$ws.Range("BF63").Value = 0.061304821323702317
$ws.Range("BG63").Value = 0
$ws.Range("BH63").Value = 45807
$ws.Range("BI63").Value = 0.46952693463407047
$ws.Range("BJ63").Value = 45828
$ws.Range("BK63").Value = 0.16240492356081113
$ws.Range("BL63").Value = 45838
$ws.Range("BM63").Value = 0.093541606810724168
$ws.Range("BN63").Value = 19.333333333333332
$ws.Range("BO63").Value = 605
$ws.Range("BP63").Value = 156713150
$ws.Range("BQ63").Value = 14263
$ws.Range("BR63").Value = 1312
$ws.Range("BS63").Value = 15575
$ws.Range("BT63").Value = 0.04999143847233023
$ws.Range("BU63").Value = 0.019337202664804092
$ws.Range("BV63").Value = 45807
$ws.Range("BW63").Value = 0.13441376791811294
$ws.Range("BX63").Value = 45814
$ws.Range("BY63").Value = 0.057274142252886927
$ws.Range("BZ63").Value = 45828
$ws.Range("CA63").Value = 0.42467390186741927
$ws.Range("CB63").Value = 11.333333333333334
$ws.Range("CC63").Value = 570
$ws.Range("CD63").Value = -103632840
$ws.Range("CE63").Value = 0.069120865940726264
$ws.Range("CF63").Value = -358
$ws.Range("CG63").Value = 6494
$ws.Range("CH63").Value = 6852
$ws.Range("CI63").Value = 0
$ws.Range("CJ63").Value = 45807
$ws.Range("CK63").Value = 0.31012607043050849
$ws.Range("CL63").Value = 45828
$ws.Range("CM63").Value = 0.18781618697601696
$ws.Range("CN63").Value = 45856
$ws.Range("CO63").Value = 0.10914572232156496
$ws.Range("CP63").Value = 25.333333333333332
$ws.Range("CQ63").Value = 575
$ws.Range("CR63").Value = -79250525
$ws.Range("CS63").Value = 0.052858388463127862
$ws.Range("CT63").Value = -603
$ws.Range("CU63").Value = 8603
$ws.Range("CV63").Value = 9206
$ws.Range("CW63").Value = 0.030362700563036974
$ws.Range("CX63").Value = 45807
$ws.Range("CY63").Value = 0.21404973661150448
$ws.Range("CZ63").Value = 45828
$ws.Range("DA63").Value = 0.34567347884330651
$ws.Range("DB63").Value = 45856
$ws.Range("DC63").Value = 0.099891220047742313
$ws.Range("DD63").Value = 25.333333333333332
$ws.Range("DE63").Value = 560
$ws.Range("DF63").Value = -58082640
$ws.Range("DG63").Value = 0.038739866367875908
$ws.Range("DH63").Value = -229
$ws.Range("DI63").Value = 4135
$ws.Range("DJ63").Value = 4364
$ws.Range("DK63").Value = 0
$ws.Range("DL63").Value = 45828
$ws.Range("DM63").Value = 0.46188518300120357
$ws.Range("DN63").Value = 45856
$ws.Range("DO63").Value = 0.15899415979385606
$ws.Range("DP63").Value = 45884
$ws.Range("DQ63").Value = 0.13135213292729225
$ws.Range("DR63").Value = 51
$ws.Range("DS63").Value = 555
$ws.Range("DT63").Value = -56674380
$ws.Range("DU63").Value = 0.037800587364524386
$ws.Range("DV63").Value = -97
$ws.Range("DW63").Value = 4436
$ws.Range("DX63").Value = 4533
$ws.Range("DY63").Value = 0
$ws.Range("DZ63").Value = 45828
$ws.Range("EA63").Value = 0.61169923231413181
$ws.Range("EB63").Value = 45838
$ws.Range("EC63").Value = 0.075398389548961095
$ws.Range("ED63").Value = 45856
$ws.Range("EE63").Value = 0.22935788426525516
$ws.Range("EF63").Value = 35.666666666666664
$ws.Range("EG63").Value = 565
$ws.Range("EH63").Value = -52039890
$ws.Range("EI63").Value = 0.034709482633691613
$ws.Range("EJ63").Value = -274
$ws.Range("EK63").Value = 4450
$ws.Range("EL63").Value = 4724
$ws.Range("EM63").Value = 0
$ws.Range("EN63").Value = 45814
$ws.Range("EO63").Value = 0.12085216159708639
$ws.Range("EP63").Value = 45828
$ws.Range("EQ63").Value = 0.19309368044783165
$ws.Range("ER63").Value = 45856
$ws.Range("ES63").Value = 0.46905139272947999
$ws.Range("ET63").Value = 27.666666666666668
$ws.Range("EU63").Value = 595
$ws.Range("EV63").Value = 430948385
$ws.Range("EW63").Value = 12427
$ws.Range("EX63").Value = 1102
$ws.Range("EY63").Value = 13529
$ws.Range("EZ63").Value = 0.092995081252691625
$ws.Range("FA63").Value = 377614965
$ws.Range("FB63").Value = 0.12045903798774152
$ws.Range("FC63").Value = 0.16051600338455865
$ws.Range("FD63").Value = 45806
$ws.Range("FE63").Value = 0.16051600338455865
$ws.Range("FF63").Value = 45807
$ws.Range("FG63").Value = 0.4472454766192861
$ws.Range("FH63").Value = 45828
$ws.Range("FI63").Value = 0.18321050914918058
$ws.Range("FJ63").Value = 8.6666666666666661
$ws.Range("FK63").Value = -53333420
$ws.Range("FL63").Value = 0.035572239205067128
$ws.Range("FM63").Value = 0.010665357668793788
$ws.Range("FN63").Value = 45828
$ws.Range("FO63").Value = 0.55242313356240791
$ws.Range("FP63").Value = 45884
$ws.Range("FQ63").Value = 0.10211299031639073
$ws.Range("FR63").Value = 45919
$ws.Range("FS63").Value = 0.074791378464010003
$ws.Range("FT63").Value = 72
$ws.Range("FU63").Value = 600
$ws.Range("FV63").Value = 385264200
$ws.Range("FW63").Value = 16594
$ws.Range("FX63").Value = 331
$ws.Range("FY63").Value = 16925
$ws.Range("FZ63").Value = 0.083136813664479184
$ws.Range("GA63").Value = 355608000
$ws.Range("GB63").Value = 0.11343882406976319
$ws.Range("GC63").Value = 0.053536478369440504
$ws.Range("GD63").Value = 45807
$ws.Range("GE63").Value = 0.14306539785381656
$ws.Range("GF63").Value = 45814
$ws.Range("GG63").Value = 0.059495849362219071
$ws.Range("GH63").Value = 45828
$ws.Range("GI63").Value = 0.36221232368225686
$ws.Range("GJ63").Value = 11.333333333333334
$ws.Range("GK63").Value = -29656200
$ws.Range("GL63").Value = 0.019780044863301694
$ws.Range("GM63").Value = 0.0058672385538268554
$ws.Range("GN63").Value = 45884
$ws.Range("GO63").Value = 0.18837882129200639
$ws.Range("GP63").Value = 45919
$ws.Range("GQ63").Value = 0.26006029093410482
$ws.Range("GR63").Value = 46010
$ws.Range("GS63").Value = 0.12349525562951423
$ws.Range("GT63").Value = 132.66666666666666
$ws.Range("GU63").Value = 580
$ws.Range("GV63").Value = 356682020
$ws.Range("GW63").Value = -1697
$ws.Range("GX63").Value = 12400
$ws.Range("GY63").Value = 14097
$ws.Range("GZ63").Value = 0.076969016675336135
$ws.Range("HA63").Value = 201439220
$ws.Range("HB63").Value = 0.064259038712093999
$ws.Range("HC63").Value = 0
$ws.Range("HD63").Value = 45821
$ws.Range("HE63").Value = 0.040977918798533873
$ws.Range("HF63").Value = 45828
$ws.Range("HG63").Value = 0.81419715584681074
$ws.Range("HH63").Value = 45919
$ws.Range("HI63").Value = 0.044358194000155483
$ws.Range("HJ63").Value = 51
$ws.Range("HK63").Value = -155242800
$ws.Range("HL63").Value = 0.10354359455036627
$ws.Range("HM63").Value = 0.067772547261451099
$ws.Range("HN63").Value = 45807
$ws.Range("HO63").Value = 0.31281476500037358
$ws.Range("HP63").Value = 45814
$ws.Range("HQ63").Value = 0.15420309347679892
$ws.Range("HR63").Value = 45828
$ws.Range("HS63").Value = 0.13238810431143988
$ws.Range("HT63").Value = 11.333333333333334
$ws.Range("HU63").Value = 590
$ws.Range("HV63").Value = 127860
$ws.Range("HW63").Value = 220957
$ws.Range("HX63").Value = 3134799773.5
$ws.Range("HY63").Value = -1499298925
$ws.Range("HZ63").Value = 1635500848.5
$ws.Range("IA63").Value = 2.0908437411839005
$ws.Range("IB63").Value = 4634098698.5
$ws.Range("IC63").Value = 0.11875309759330971
$ws.Range("ID63").Value = 45806
$ws.Range("IE63").Value = 0.11875309759330971
$ws.Range("IF63").Value = 45807
$ws.Range("IG63").Value = 0.19968789363496547
$ws.Range("IH63").Value = 45828
$ws.Range("II63").Value = 0.23165393463619602

# --- Update active selection to reflect F69 as last active cell ---
$ws.Range("F69").Select()
